# The sheet gained a third column "type" classifying each row of the
# string table (title / descriptions / button / button).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "type"
$ws.Range("C2").Value = "title"
$ws.Range("C3").Value = "descriptions"
$ws.Range("C4").Value = "button"
$ws.Range("C5").Value = "button"
